$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.302.16'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.23%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.790.71'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.06%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '226.14'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.39%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.556'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.49%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '32.30'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.16%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.91%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.43%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.58%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.049.80'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.11%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.789.63'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.42%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '11.00'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.08%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.97%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '34.280.35'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.43%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.01%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.97'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.00%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0802'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.77%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '246.52'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.35%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.94'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.16'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.07'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.52%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '162.37'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.47%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.17'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.12%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.37'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.35%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.43%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.21%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.56%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.90'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +7.95%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0521'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.09%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.77'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.68%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.37%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.439.65'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.20%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.60'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +7.25%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.660'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.88%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.06'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.77%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.91%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '82.26'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +2.15%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.82%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '14.15'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +6.23%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.68%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.12%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.97%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.17%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.44%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.945.07'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.20%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '105.70'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.93%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₆0130'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -6.14%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.05%  '
